$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 79: fill in the already-present placeholder row with the 2ème Série IDF Place result
$ws.Range("A79").Value = 538.2
$ws.Range("B79").Value = 45725
$ws.Range("C79").Value = "Open"
$ws.Range("D79").Value = "2ème Série IDF Place"
$ws.Range("E79").Value = 4
$ws.Range("F79").Value = 2

# Row 80: fill in the already-present placeholder row with a 2ème Série IDF result
$ws.Range("A80").Value = 464
$ws.Range("B80").Value = 45725
$ws.Range("C80").Value = "Open"
$ws.Range("D80").Value = "2ème Série IDF"

# Row 81: brand-new row with another 2ème Série IDF result
$ws.Range("A81").Value = 598
$ws.Range("B81").Value = 45725
$ws.Range("C81").Value = "Open"
$ws.Range("D81").Value = "2ème Série IDF"
